# CA_variable_check_accepted_ranges.xlsx
# "updated after review of variables"
#
# Delvin/reviewer went through every row of the missingness table and added a
# new "reviewer decision" column (E), plus filled in a few still-missing
# accepted-range / unit notes in columns F and G for the identifier/outcome
# variables near the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give column E (the new notes column) a sensible width, matching the
# "best fit" width Excel would compute for the longest note in the column.
$ws.Range("E1").ColumnWidth = 36.6

# --- Column E: reviewer decision / comment for each variable row (2-64) ---
$ws.Range("E2").Value = "Too many missing"
$ws.Range("E3:E6").Value = "A lot of missing but very important variable"
$ws.Range("E7:E13").Value = "Too many missing (remind me its >25% right?)"
$ws.Range("E14").Value = "Cutoff?"
$ws.Range("E15:E64").Value = "Keep"

# --- Columns F/G: accepted ranges / units filled in for a few variables
#     that didn't have them recorded yet ---
$ws.Range("F49").Value = "30 to 210 "
$ws.Range("G49").Value = "kg"
$ws.Range("F51").Value = "no cutoff (identifier)"
$ws.Range("F52").Value = "2008 to 2018"
$ws.Range("F56").Value = "0 or 1"
$ws.Range("F58").Value = "no cutoff"
$ws.Range("F60").Value = "no cutoff"
$ws.Range("F61").Value = "18 to 90"

# Leave the selection where the reviewer ended up looking (row 12, column F)
$ws.Range("F12").Select() | Out-Null
